# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" fund-holding detail sheet (just before the
# existing "2022-Q2" sheet) and updates the "总计" (totals) roll-up sheet
# so it lists the new quarter first, pushing the older quarters down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: set a cell to a literal TEXT value even when the text looks
# like a number (e.g. "1.30", "012430"), bypassing Excel's automatic
# "looks like a number -> store as number" coercion that a plain
# `Range.Value = "1.30"` assignment would trigger. We do this by writing
# a self-quoting formula into an unused scratch cell and then Copy()-ing
# that cell onto the real target: copying a no-reference text formula
# bakes down to a plain text value on the destination (no formula, no
# extra number-format/style) rather than carrying the formula itself.
# ---------------------------------------------------------------------
function Set-TextValue {
    param($ws, [string]$addr, [string]$text)
    $scratch = $ws.Range("ZZ1")
    $escaped = $text.Replace('"', '""')
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy($ws.Range($addr))
    $scratch.Value = ""
}

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by duplicating "2022-Q2" (so it
#    inherits identical column headers/formatting/page setup) and
#    placing the copy immediately before it.
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item("2022-Q2")
$q2Sheet.Copy($q2Sheet)
$newSheet = $wb.Worksheets.Item("2022-Q2 (2)")
$newSheet.Name = "2022-Q3"

# "2022-Q2" has 3 fund rows (rows 2-4); "2022-Q3" only needs 1 (row 2),
# so drop the extra copied rows.
$newSheet.Range("A3:H4").Clear()

# Fill in 2022-Q3's single fund-holding row.
$newSheet.Range("A2").Value = 0
Set-TextValue $newSheet "B2" "012430"
$newSheet.Range("C2").Value = "农银汇理瑞康6个月持有期混合"
Set-TextValue $newSheet "D2" "1.30"
Set-TextValue $newSheet "E2" "24.44"
Set-TextValue $newSheet "F2" "1.13"
Set-TextValue $newSheet "G2" "0.0147"
$newSheet.Range("H2").Value = 4

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: shift the quarter list down one
#    row and add the new 2022-Q3 totals at the top.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("总计")

# Grow the table by one row, copying row 5's formatting onto new row 6
# (so A6 picks up the same index-column style as A2:A5).
$ws1.Range("A5:D5").Copy($ws1.Range("A6:D6"))

$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 1
$ws1.Range("D2").Value = 0.01

$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = "2022-Q2"
$ws1.Range("C3").Value = 3
$ws1.Range("D3").Value = 0.04

$ws1.Range("A4").Value = 2
$ws1.Range("B4").Value = "2022-Q1"
$ws1.Range("C4").Value = 2
$ws1.Range("D4").Value = 0.32

$ws1.Range("A5").Value = 3
$ws1.Range("B5").Value = "2021-Q4"
$ws1.Range("C5").Value = 3
$ws1.Range("D5").Value = 0.54

$ws1.Range("A6").Value = 4
$ws1.Range("B6").Value = "2021-Q3"
$ws1.Range("C6").Value = 3
$ws1.Range("D6").Value = 0.58

# ---------------------------------------------------------------------
# 3. Restore "2021-Q3" (still the last tab) as the selected sheet, since
#    creating/copying sheets above shifted the active tab.
# ---------------------------------------------------------------------
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()

Write-Output "done"
